# Adds 10 new "Kinect" papers (rows 3-13) to the catalogue sheet, matching
# the commit "2016/1/14: added 10 papers on kinect".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows (row, A..I). Empty string means "leave blank".
# F_IsText marks Time values that must stay literal text (e.g. "2015.6.2")
# instead of being auto-parsed into a date serial by Excel.
$rows = @(
    @{ A=1;  C=1;  D="Assessing the Suitability of the Microsoft Kinect for Calculating Person Specific Body Segment Parameters";
             E="Sean Clarkson, Jon Wheat, Ben Heller, Simon Choppin";
             F=2014; FIsText=$false; H="B"; I="ECCV" },
    @{ A=2;  C=2;  D="A framework for gait-based recognition using Kinect";
             E="Dimitris Kastaniotis, Ilias Theodorakopoulosa, Christos Theoharatosb, George Economoua, Spiros Fotopoulos";
             F="2015.6.2"; FIsText=$true; H="C"; I="PRL" },
    @{ A=3;  C=3;  D="Real Time Gait Recognition System based on Kinect Skeleton Feature";
             E="Shuming Jiang, Yufei Wang, Yuanyuan Zhang, and Jiande Sun";
             F=2014; FIsText=$false; H="C"; I="ACCV" },
    @{ A=4;  C=4;  D="Detection of gait cycles in treadmill walking using a Kinect";
             E="Edouard Auvinet, Franck Multon, Carl-Eric Aubin, Jean Meunier, Maxime Raison";
             F="2014.8.11"; FIsText=$true; H=""; I="Gait & posture" },
    @{ A=5;  C=5;  D="Person Identification in Natural Static Postures Using Kinect";
             E="Reddy V R, Chakravarty K, Aniruddha S";
             F=2014; FIsText=$false; H="B"; I="ECCV" },
    @{ A=6;  C=6;  D="Fall detection in homes of older adults using the microsoft kinect";
             E="Erik E. Stone, and Marjorie Skubic";
             F=2015; FIsText=$false; H=""; I="Biomedical and Health Informatics, IEEE Journal of" },
    @{ A=7;  C=7;  D="Full body gait analysis with Kinect";
             E="Gabel M, Gilad-Bachrach R, Renshaw E";
             F=2012; FIsText=$false; H=""; I="EMBC" },
    @{ A=8;  C=8;  D="Instrumenting gait assessment using the Kinect in people living with stroke: reliability and association with balance tests";
             E="Clark R A, Vernon S, Mentiplay B F";
             F=2015; FIsText=$false; H=""; I="Journal of neuroengineering and rehabilitation" },
    @{ A=9;  C=9;  D="Person Identification Using Full-Body Motion and Anthropometric Biometrics from Kinect Videos";
             E="Munsell B C, Temlyakov A, Qu C";
             F=2012; FIsText=$false; H="B"; I="ECCV" },
    @{ A=10; C=10; D="Reachable workspace in facioscapulohumeral muscular dystrophy (FSHD) by kinect";
             E="Han J J, Kurillo G, Abresch R T";
             F=2015; FIsText=$false; H=""; I="Muscle & nerve" },
    @{ A=11; C=11; D="Towards skeleton biometric identification using the microsoft kinect sensor";
             E="Araujo R M, Graña G, Andersson V";
             F=2013; FIsText=$false; H=""; I="ACM Symposium on Applied Computing" }
)

$r = 3
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = "Kinect"
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E

    if ($row.FIsText) {
        # Force text storage so "2015.6.2" / "2014.8.11" aren't coerced
        # into date serials; ClearFormats keeps the cell on the default
        # style once the literal value has been committed.
        $ws.Cells.Item($r, 6).NumberFormat = "@"
        $ws.Cells.Item($r, 6).Value = $row.F
        $ws.Cells.Item($r, 6).ClearFormats()
    } else {
        $ws.Cells.Item($r, 6).Value = $row.F
    }

    if ($row.H -ne "") {
        $ws.Cells.Item($r, 8).Value = $row.H
    }
    $ws.Cells.Item($r, 9).Value = $row.I
    $r = $r + 1
}

# Match the post-edit column widths recorded in the diff (57.625 / 33 /
# 15.5 "characters"). The host only round-trips ColumnWidth to the nearest
# 1/7 of a character (5px padding + 7px/char), so feed it the pre-image
# that lands closest to each target after that quantization.
$ws.Columns.Item(4).ColumnWidth = 56.857142857142854   # -> stored 57.571... (target 57.625)
$ws.Columns.Item(5).ColumnWidth = 32.285714285714285   # -> stored 33 (exact)
$ws.Columns.Item(9).ColumnWidth = 14.714285714285714   # -> stored 15.428... (target 15.5)

# Match the recorded selection after the edit.
$ws.Range("D18").Select()
